# Updates the "想去人数" (want-to-go count) figures in column F of the
# 展览 / 本地生活 / 全部类型 sheets to reflect a newer scrape of the source
# data (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 735
$ws.Cells.Item(3, 6).Value = 14042
$ws.Cells.Item(4, 6).Value = 14043
$ws.Cells.Item(5, 6).Value = 14066
$ws.Cells.Item(6, 6).Value = 1360
$ws.Cells.Item(7, 6).Value = 1384
$ws.Cells.Item(8, 6).Value = 5821
$ws.Cells.Item(9, 6).Value = 973
$ws.Cells.Item(15, 6).Value = 423
$ws.Cells.Item(17, 6).Value = 1180
$ws.Cells.Item(21, 6).Value = 2252
$ws.Cells.Item(24, 6).Value = 3272
$ws.Cells.Item(27, 6).Value = 2331
$ws.Cells.Item(31, 6).Value = 1768
$ws.Cells.Item(32, 6).Value = 1063
$ws.Cells.Item(33, 6).Value = 1346
$ws.Cells.Item(35, 6).Value = 136
$ws.Cells.Item(36, 6).Value = 4673
$ws.Cells.Item(37, 6).Value = 4747
$ws.Cells.Item(42, 6).Value = 3264
$ws.Cells.Item(45, 6).Value = 332
$ws.Cells.Item(46, 6).Value = 88
$ws.Cells.Item(48, 6).Value = 4408
$ws.Cells.Item(49, 6).Value = 542

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 7421
$ws.Cells.Item(3, 6).Value = 217
$ws.Cells.Item(4, 6).Value = 675

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 7421
$ws.Cells.Item(3, 6).Value = 735
$ws.Cells.Item(4, 6).Value = 217
$ws.Cells.Item(5, 6).Value = 675
$ws.Cells.Item(7, 6).Value = 14043
$ws.Cells.Item(8, 6).Value = 14043
$ws.Cells.Item(9, 6).Value = 14067
$ws.Cells.Item(10, 6).Value = 1361
$ws.Cells.Item(11, 6).Value = 1384
$ws.Cells.Item(12, 6).Value = 5821
$ws.Cells.Item(13, 6).Value = 973
$ws.Cells.Item(16, 6).Value = 423
$ws.Cells.Item(17, 6).Value = 1180
$ws.Cells.Item(21, 6).Value = 3272
$ws.Cells.Item(26, 6).Value = 1768
$ws.Cells.Item(32, 6).Value = 1063
$ws.Cells.Item(33, 6).Value = 1346
$ws.Cells.Item(36, 6).Value = 4673
$ws.Cells.Item(37, 6).Value = 4747
$ws.Cells.Item(40, 6).Value = 3264
$ws.Cells.Item(43, 6).Value = 332
$ws.Cells.Item(44, 6).Value = 88
$ws.Cells.Item(46, 6).Value = 4408
